$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.756.69'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '3.309.43'
$ws.Range("E3").Value = '  +5.14%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.49'
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.20'
$ws.Range("E6").Value = '  +2.57%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.306.75'
$ws.Range("E8").Value = '  +5.12%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("E13").Value = '  +0.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.65'
$ws.Range("E14").Value = '  +1.16%  '
$ws.Range("D15").Value = '3.854.41'
$ws.Range("E15").Value = '  +5.05%  '
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '3.313.21'
$ws.Range("E17").Value = '  +5.15%  '
$ws.Range("D18").Value = '63.870.29'
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("E19").Value = '  +2.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.98'
$ws.Range("E20").Value = '  +1.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.16'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  +4.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.10'
$ws.Range("E23").Value = '  +4.88%  '
$ws.Range("E24").Value = '  +5.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.47'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.76'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("E28").Value = '  +4.51%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.11'
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.16'
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.83'
$ws.Range("E32").Value = '  +7.28%  '
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("E35").Value = '  +3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.02'
$ws.Range("E36").Value = '  +3.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.18'
$ws.Range("E37").Value = '  +1.20%  '
$ws.Range("D38").Value = '0.0₃0739'
$ws.Range("E38").Value = '  +4.98%  '
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '434.68'
$ws.Range("E40").Value = '  +2.95%  '
$ws.Range("D41").Value = '3.065.59'
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.115'
$ws.Range("E44").Value = '  +2.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.264'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.34'
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '36.20'
$ws.Range("E48").Value = '  +12.27%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '125.46'
$ws.Range("E50").Value = '  +4.24%  '
$ws.Range("E51").Value = '  +0.82%  '
